# Stock_reporter_no_panda mixing BTC/USD and UL quantity
#
# The BTC-USD row was removed from its original position (row 6) and a
# re-labelled copy of its data ("xxx") was appended at the bottom of the
# table (row 25). Every other row shifts up by one. In addition, the
# "Quantity" column (D) is reformatted from 2-decimal text to 6-decimal
# text across the whole table, and UL's quantity picks up the extra
# floating-point precision (0.03 -> 0.029997) that the 6-decimal
# formatting exposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the BTC-USD row (row 6); everything below shifts up one row.
$ws.Rows.Item(6).Delete()

# 2) Re-append the former BTC-USD data as a new "xxx" row at the bottom
#    (now row 25, since the table lost a row).
$ws.Range("A25").Value = "xxx"
$ws.Range("B25").Value = 44833
$ws.Range("C25").Value = 47649.16
$ws.Range("D25").Value = "'2.000000"
$ws.Range("E25").Value = "'2816.16"
$ws.Range("F25").Value = "'6.28 %"

# 3) Reformat the Quantity column (D) to 6 decimal places for every data
#    row, also fixing up UL's value to its full-precision text. The
#    leading apostrophe forces Excel to keep these as text instead of
#    auto-converting the numeric-looking strings to numbers.
$quantities = @{
    2  = "2.000000"
    3  = "1.000000"
    4  = "1.000000"
    5  = "6.000000"
    6  = "1.000000"
    7  = "3.000000"
    8  = "3.000000"
    9  = "1.000000"
    10 = "2.000000"
    11 = "4.000000"
    12 = "5.000000"
    13 = "1.000000"
    14 = "2.000000"
    15 = "1.000000"
    16 = "1.000000"
    17 = "3.000000"
    18 = "3.000000"
    19 = "2.000000"
    20 = "1.000000"
    21 = "1.000000"
    22 = "0.029997"
    23 = "2.000000"
    24 = "2.000000"
    25 = "2.000000"
}

foreach ($row in $quantities.Keys) {
    $ws.Range("D$row").Value = "'" + $quantities[$row]
}
